$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.607.70'
$ws.Range('E2').Value = '  -7.27%  '
$ws.Range('D3').Value = '1.698.33'
$ws.Range('E3').Value = '  -5.88%  '
$c = $ws.Range('D4')
$c.Value = "'" + '1.005'
$c.Style = "Normal"
$ws.Range('E4').Value = '  +0.23%  '
$c = $ws.Range('D5')
$c.Value = "'" + '220.14'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -5.24%  '
$c = $ws.Range('D6')
$c.Value = "'" + '0.5137'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -13.22%  '
$ws.Range('E7').Value = '  +0.18%  '
$c = $ws.Range('D8')
$c.Value = "'" + '0.2657'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -4.20%  '
$c = $ws.Range('D9')
$c.Value = "'" + '22.21'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -4.75%  '
$c = $ws.Range('D10')
$c.Value = "'" + '0.06277'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -7.89%  '
$c = $ws.Range('D11')
$c.Value = "'" + '0.07354'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -2.02%  '
$ws.Range('D12').Value = '1.701.13'
$ws.Range('E12').Value = '  -5.62%  '
$c = $ws.Range('D13')
$c.Value = "'" + '4.527'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -4.87%  '
$c = $ws.Range('D14')
$c.Value = "'" + '0.5853'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -5.87%  '
$ws.Range('D15').Value = '1.927.64'
$ws.Range('E15').Value = '  -5.96%  '
$c = $ws.Range('D16')
$c.Value = "'" + '0.000008469'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -8.07%  '
$c = $ws.Range('D17')
$c.Value = "'" + '65.66'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -13.13%  '
$ws.Range('D18').Value = '26.634.24'
$ws.Range('E18').Value = '  -7.09%  '
$c = $ws.Range('D19')
$c.Value = "'" + '5.029'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -8.16%  '
$ws.Range('E20').Value = '  +0.19%  '
$c = $ws.Range('D21')
$c.Value = "'" + '10.97'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -4.62%  '
$c = $ws.Range('D22')
$c.Value = "'" + '187.27'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -11.13%  '
$c = $ws.Range('D23')
$c.Value = "'" + '6.280'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +0.28%  '
$c = $ws.Range('D25')
$c.Value = "'" + '145.05'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -5.68%  '
$c = $ws.Range('D26')
$c.Value = "'" + '7.569'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -3.67%  '
$c = $ws.Range('D27')
$c.Value = "'" + '0.1151'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -9.04%  '
$c = $ws.Range('D28')
$c.Value = "'" + '15.73'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -4.20%  '
$c = $ws.Range('D29')
$c.Value = "'" + '1.330'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -6.85%  '
$c = $ws.Range('D30')
$c.Value = "'" + '0.05706'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -7.54%  '
$c = $ws.Range('D31')
$c.Value = "'" + '1.336'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -6.44%  '
$c = $ws.Range('D32')
$c.Value = "'" + '3.523'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -6.82%  '
$ws.Range('E33').Value = '  -6.15%  '
$c = $ws.Range('D34')
$c.Value = "'" + '1.651'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -4.92%  '
$c = $ws.Range('D35')
$c.Value = "'" + '1.029'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -2.97%  '
$c = $ws.Range('D36')
$c.Value = "'" + '0.6032'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -6.14%  '
$c = $ws.Range('D37')
$c.Value = "'" + '2.370'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -5.16%  '
$ws.Range('E38').Value = '  -1.17%  '
$ws.Range('D39').Value = '1.102.16'
$ws.Range('E39').Value = '  -3.94%  '
$c = $ws.Range('D40')
$c.Value = "'" + '0.01609'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -4.98%  '
$c = $ws.Range('D41')
$c.Value = "'" + '0.8626'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -2.41%  '
$c = $ws.Range('D42')
$c.Value = "'" + '5.862'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -10.51%  '
$ws.Range('E43').Value = '  -0.04%  '
$c = $ws.Range('D44')
$c.Value = "'" + '99.03'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -0.94%  '
$ws.Range('D45').Value = '1.855.01'
$ws.Range('E45').Value = '  -5.31%  '
$c = $ws.Range('D46')
$c.Value = "'" + '0.00000000111'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +0.27%  '
$c = $ws.Range('D47')
$c.Value = "'" + '56.71'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -6.36%  '
$c = $ws.Range('D48')
$c.Value = "'" + '8.181'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -2.06%  '
$c = $ws.Range('D49')
$c.Value = "'" + '1.001'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -0.07%  '
$c = $ws.Range('D50')
$c.Value = "'" + '0.05244'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -4.18%  '
$c = $ws.Range('D51')
$c.Value = "'" + '0.4329'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -3.28%  '
